$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A34").Value = "погода"
$ws.Range("B34").Value = "Погода"
$ws.Range("C34").Value = 1
